$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 6 (ano 2025) metrics per diff
$ws.Range("C6").Value = 413
$ws.Range("D6").Value = 308
$ws.Range("F6").Value = 67.39606126914661
$ws.Range("G6").Value = 25.42372881355932
$ws.Range("H6").Value = 74.57627118644068
